$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-10 Thursday", "2025-04-11 Friday"),
    @("32×36=1152", "27×43=1161"),
    @("11×58=638", "49×57=2793"),
    @("32×67=2144", "98×72=7056"),
    @("95×46=4370", "53×45=2385"),
    @("70×90=6300", "80×55=4400"),
    @("74×34=2516", "60×84=5040"),
    @("52×62=3224", "89×48=4272"),
    @("11×56=616", "33×99=3267"),
    @("31×53=1643", "65×27=1755"),
    @("86×44=3784", "96×21=2016"),
    @("21×18=378", "53×77=4081"),
    @("73×62=4526", "13×72=936"),
    @("59×20=1180", "14×75=1050"),
    @("89×63=5607", "56×56=3136"),
    @("21×60=1260", "46×26=1196"),
    @("43×18=774", "28×38=1064"),
    @("74×68=5032", "52×37=1924"),
    @("44×58=2552", "30×79=2370"),
    @("95×67=6365", "98×96=9408"),
    @("86×55=4730", "53×39=2067"),
    @("32×73=2336", "99×87=8613"),
    @("83×79=6557", "53×72=3816"),
    @("39×65=2535", "39×88=3432"),
    @("37×91=3367", "98×39=3822"),
    @("22×44=968", "91×87=7917")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
